$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 5.586269137925634

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 6.189590430959694

# Row 4
$ws.Range("B4").Value = 1.455362044514542
$ws.Range("C4").Value = 117.745847958593
$ws.Range("D4").Value = 261.3203778131603
$ws.Range("E4").Value = 2195978.878461985
$ws.Range("G4").Value = 2196359.400049801
